$wb = $excel.ActiveWorkbook

# The "as_of_utc" column (AA) on both data sheets needs its timestamp
# refreshed from 2025-11-01 03:02:28 to 2025-11-01 07:02:23 for every
# data row (rows 2-26).
$oldValue = "2025-11-01 03:02:28"
$newValue = "2025-11-01 07:02:23"

foreach ($sheetName in @("Главные", "Линейные")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Range("AA$row")
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
